# chore: clean format, unused variables and comments
# Appends the newest sell-data record (row 20) to the SellData sheet and
# refreshes the timestamp on the previous last row, mirroring the source
# data regeneration that produced this state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The upstream data refresh also re-stamped the prior last row's date.
$ws.Cells.Item(19, 5).Value = 45819.82946606482

# New sale record.
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "1AYB-2AYB-4AYB-1P-2P-4P-1AP-5AYB"
$ws.Cells.Item(20, 3).Value = "2-1-1-1-1-1-1-1"
$ws.Cells.Item(20, 4).Value = 59150
$ws.Cells.Item(20, 5).Value = 45821.69294795885
$ws.Cells.Item(20, 5).NumberFormat = $ws.Cells.Item(19, 5).NumberFormat

Write-Output "Appended row 20 to SellData and updated E19 timestamp"
